# Daily update at 8 AM UTC
# - Row 14 (previously the "latest" row) reverts to the normal date format.
# - A new row 15 is appended with the next day's data and takes on the
#   "latest" row date format.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 14 was styled as the most-recent entry (date-only format); now that a
# newer row is being added, it goes back to the regular timestamp format
# used by every other historical row.
$ws.Range("A14").NumberFormat = $ws.Range("A13").NumberFormat

# Append the new day's results.
$ws.Range("A15").Value = 45755
$ws.Range("B15").Value = 56
$ws.Range("C15").Value = 59
$ws.Range("D15").Value = 54

# The newest row takes on the "latest" row's date-only format.
$ws.Range("A15").NumberFormat = "YYYY-MM-DD"
